# issue #5: stock data from json to db
# Adds a "category" column (value "normal") right after "property_category",
# and appends "source_file" / "index" columns at the end of the 股票 (stock)
# sheet. Also fixes a garbled numeric value in the "total" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- 1. Insert a new column I ("category") before the existing "date"
#        column (old column I), shifting date/legislator_name/legislator_id
#        one column to the right (I->J, J->K, K->L). Insert() carries the
#        surrounding cell formatting (header bold/border style, data style)
#        along with the shift, same as Excel's own "Insert Cut Cells".
$ws.Columns.Item(9).Insert()

# --- 2. Append two more columns (source_file, index) after legislator_id.
#        Build them the same way: copy the neighbouring legislator_id column
#        and Insert-shift it outward so the new columns inherit matching
#        cell styles, then overwrite with the real values below.
$ws.Range("L1:L7").Copy()
$ws.Range("M1:M7").Insert(-4161)
$ws.Range("L1:L7").Copy()
$ws.Range("N1:N7").Insert(-4161)

# --- 3. Header row
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- 4. Data rows: category = "normal" for every stock record; source_file
#        is the originating filename stem; index is the original row key
#        (same number already present in column A).
$rows = 2..7
foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmpf421"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value2
}

# --- 5. Data fix: the "total" value for row 7 was garbled ("143《290")
#        and is corrected to the text "143290". Go through a TEXT() formula
#        + paste-values round-trip so the corrected value stays a text cell
#        (matching the rest of this mixed string/number column) instead of
#        Excel auto-converting the plain digit string to a number.
$ws.Range("G7").Formula = '=TEXT(143290,"0")'
$ws.Range("G7").Copy()
$ws.Range("G7").PasteSpecial(-4163)
